$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F
$ws.Columns.Item(6).ColumnWidth = 33.85546875

# Add a yellow-fill style to rows that remain "open" (keep value "open")
$ws.Range("B7").Interior.Color = 65535
$ws.Range("B8").Interior.Color = 65535
$ws.Range("B11").Interior.Color = 65535
$ws.Range("B14").Interior.Color = 65535
$ws.Range("B16").Interior.Color = 65535

# Rows that flip from "open" to "closed"
$ws.Range("B12").Value = "closed"
$ws.Range("B13").Value = "closed"
$ws.Range("B15").Value = "closed"
$ws.Range("B17").Value = "closed"

# New row 18 data
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "open"
$ws.Range("B18").Interior.Color = 65535
$ws.Range("C18").Value = "refund"
$ws.Range("D18").Value = "checking refund amount"
$ws.Range("E18").Value = "1.cancel any applicant who has paid more than 1000"
$ws.Range("F18").Value = "calculate amt of refund based on dates"
$ws.Range("G18").Value = "1000*% based on time + extra amount paid as it is"
$ws.Range("H18").Value = "total amt * %"

$ws.Range("C19").Select()
